$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.730.71'
$ws.Range('E2').Value = '  +0.00%  '
$ws.Range('D3').Value = '2.209.53'
$ws.Range('E3').Value = '  -2.14%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '229.91'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.41%  '
$ws.Range('E6').Value = '  -4.11%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '60.37'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -6.18%  '
$ws.Range('E8').Value = '  -0.06%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.403'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.49%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '57.26'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -4.78%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0888'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.69%  '
$ws.Range('E12').Value = '  -2.59%  '
$ws.Range('D13').Value = '2.537.04'
$ws.Range('E13').Value = '  -2.12%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '15.40'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -4.64%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '22.24'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.37%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.67'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.36%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.794'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -3.97%  '
$ws.Range('D18').Value = '2.225.77'
$ws.Range('E18').Value = '  -1.56%  '
$ws.Range('D19').Value = '41.682.71'
$ws.Range('E19').Value = '  +0.18%  '
$ws.Range('B20').Value = 'ShibaInu'
$ws.Range('C20').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D20').Value = '0.0₃0901'
$ws.Range('E20').Value = '  -4.44%  '
$ws.Range('B21').Value = 'Litecoin'
$ws.Range('C21').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '72.10'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -3.15%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.06'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.06%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '242.70'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -4.00%  '
$ws.Range('E24').Value = '  -0.16%  '
$ws.Range('E25').Value = '  -3.31%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.29'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.47%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.65'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.07%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '169.27'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.38%  '
$ws.Range('E29').Value = '  -5.26%  '
$ws.Range('E30').Value = '  +0.60%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '19.77'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.79%  '
$ws.Range('E32').Value = '  -8.68%  '
$ws.Range('E33').Value = '  -3.73%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.02'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.08%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.63'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.41%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0648'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.97%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.37'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -5.41%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.34'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -8.58%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.54'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -8.28%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.000240'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -7.56%  '
$ws.Range('E41').Value = '  -0.03%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0239'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.08%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.56'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.43%  '
$ws.Range('B44').Value = 'Cronos'
$ws.Range('C44').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0955'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -3.19%  '
$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '97.27'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -5.25%  '
$ws.Range('E46').Value = '  -3.55%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.39'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -14.45%  '
$ws.Range('D48').Value = '1.466.87'
$ws.Range('E48').Value = '  -2.94%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '16.41'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -6.85%  '
$ws.Range('E50').Value = '  -0.98%  '
$ws.Range('E51').Value = '  -5.88%  '
